# Scheduled runner update: refresh Tonberry profit-calc cells (currentAveragePrice /
# LevePriceNQ/HQ / LeveProfitNQ/HQ and related columns) across the per-job leve sheets
# with freshly pulled market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 495.53845
$ws.Range("I15").Value = 495.53845
$ws.Range("K15").Value = 1486.61535
$ws.Range("M15").Value = -1317.61535
$ws.Range("H32").Value = 3376.6667
$ws.Range("I32").Value = 4000
$ws.Range("J32").Value = 3252
$ws.Range("K32").Value = 4000
$ws.Range("L32").Value = 3252
$ws.Range("M32").Value = -3674
$ws.Range("N32").Value = -3904
$ws.Range("H51").Value = 4266.6665
$ws.Range("J51").Value = 3400
$ws.Range("L51").Value = 3400
$ws.Range("N51").Value = -4368
$ws.Range("H81").Value = 30099.334
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996
$ws.Range("H84").Value = 30099.334
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984
$ws.Range("H92").Value = 14706116
$ws.Range("I92").Value = 15625185
$ws.Range("K92").Value = 15625185
$ws.Range("M92").Value = -15623937
$ws.Range("H106").Value = 4751.25
$ws.Range("I106").Value = 4751.25
$ws.Range("K106").Value = 4751.25
$ws.Range("M106").Value = -4120.25
$ws.Range("H121").Value = 1114.8334
$ws.Range("J121").Value = 1297.8
$ws.Range("L121").Value = 3893.4
$ws.Range("N121").Value = -7387.4
$ws.Range("H138").Value = 1888.1757
$ws.Range("I138").Value = 1709.9445
$ws.Range("K138").Value = 5129.833500000001
$ws.Range("M138").Value = 10.16649999999936
$ws.Range("H139").Value = 61997.145
$ws.Range("J139").Value = 61997.145
$ws.Range("L139").Value = 61997.145
$ws.Range("N139").Value = -72277.14499999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1849.1666
$ws.Range("I45").Value = 1649.5
$ws.Range("K45").Value = 1649.5
$ws.Range("M45").Value = -1272.5
$ws.Range("H74").Value = 1108.9395
$ws.Range("I74").Value = 607.6799999999999
$ws.Range("K74").Value = 607.6799999999999
$ws.Range("M74").Value = 266.3200000000001
$ws.Range("H77").Value = 1108.9395
$ws.Range("I77").Value = 607.6799999999999
$ws.Range("K77").Value = 3038.4
$ws.Range("M77").Value = 1329.6
$ws.Range("H97").Value = 1044.9445
$ws.Range("I97").Value = 916.625
$ws.Range("J97").Value = 2071.5
$ws.Range("K97").Value = 916.625
$ws.Range("L97").Value = 2071.5
$ws.Range("M97").Value = -420.625
$ws.Range("N97").Value = -3063.5
$ws.Range("H132").Value = 1565.4
$ws.Range("I132").Value = 1519.3334
$ws.Range("K132").Value = 4558.0002
$ws.Range("M132").Value = -2028.0002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1060.3
$ws.Range("I99").Value = 796.5
$ws.Range("K99").Value = 796.5
$ws.Range("M99").Value = 701.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2269.8333
$ws.Range("I31").Value = 2045
$ws.Range("J31").Value = 2854.4
$ws.Range("K31").Value = 2045
$ws.Range("L31").Value = 2854.4
$ws.Range("M31").Value = -1750
$ws.Range("N31").Value = -3444.4
$ws.Range("H34").Value = 2269.8333
$ws.Range("I34").Value = 2045
$ws.Range("J34").Value = 2854.4
$ws.Range("K34").Value = 2045
$ws.Range("L34").Value = 2854.4
$ws.Range("M34").Value = -1843
$ws.Range("N34").Value = -3258.4
$ws.Range("H58").Value = 5437912.5
$ws.Range("I58").Value = 21739630
$ws.Range("J58").Value = 4006.8333
$ws.Range("K58").Value = 21739630
$ws.Range("L58").Value = 4006.8333
$ws.Range("M58").Value = -21739427
$ws.Range("N58").Value = -4412.8333
$ws.Range("H94").Value = 1473.25
$ws.Range("I94").Value = 1796
$ws.Range("J94").Value = 1365.6666
$ws.Range("K94").Value = 1796
$ws.Range("L94").Value = 1365.6666
$ws.Range("M94").Value = -1345
$ws.Range("N94").Value = -2267.6666
$ws.Range("H132").Value = 1736.7567
$ws.Range("I132").Value = 1323.75
$ws.Range("K132").Value = 3971.25
$ws.Range("M132").Value = -1441.25
$ws.Range("H136").Value = 5437912.5
$ws.Range("I136").Value = 21739630
$ws.Range("J136").Value = 4006.8333
$ws.Range("K136").Value = 65218890
$ws.Range("L136").Value = 12020.4999
$ws.Range("M136").Value = -65216340
$ws.Range("N136").Value = -17120.4999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 406.5238
$ws.Range("I5").Value = 388.2
$ws.Range("J5").Value = 423.18182
$ws.Range("K5").Value = 1164.6
$ws.Range("L5").Value = 1269.54546
$ws.Range("M5").Value = -1052.6
$ws.Range("N5").Value = -1493.54546
$ws.Range("H113").Value = 8526.691999999999
$ws.Range("I113").Value = 25725.25
$ws.Range("J113").Value = 882.8889
$ws.Range("K113").Value = 77175.75
$ws.Range("L113").Value = 2648.6667
$ws.Range("M113").Value = -75005.75
$ws.Range("N113").Value = -6988.6667
$ws.Range("H131").Value = 12407.049
$ws.Range("J131").Value = 12596.667
$ws.Range("L131").Value = 37790.001
$ws.Range("N131").Value = -47870.001
$ws.Range("H132").Value = 1326.125
$ws.Range("J132").Value = 1560.8
$ws.Range("L132").Value = 14047.2
$ws.Range("N132").Value = -19107.2
$ws.Range("H133").Value = 4055.9
$ws.Range("I133").Value = 1853
$ws.Range("K133").Value = 5559
$ws.Range("M133").Value = -499
$ws.Range("H135").Value = 406.5238
$ws.Range("I135").Value = 388.2
$ws.Range("J135").Value = 423.18182
$ws.Range("K135").Value = 3493.8
$ws.Range("L135").Value = 3808.63638
$ws.Range("M135").Value = -958.7999999999997
$ws.Range("N135").Value = -8878.63638
$ws.Range("H140").Value = 3167.9473
$ws.Range("I140").Value = 1559.4
$ws.Range("K140").Value = 4678.200000000001
$ws.Range("M140").Value = 501.7999999999993
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 37218
$ws.Range("J127").Value = 37218
$ws.Range("L127").Value = 37218
$ws.Range("N127").Value = -47138
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1576.7693
$ws.Range("I82").Value = 1183
$ws.Range("J82").Value = 3742.5
$ws.Range("K82").Value = 1183
$ws.Range("L82").Value = 3742.5
$ws.Range("M82").Value = -822
$ws.Range("N82").Value = -4464.5
$ws.Range("H85").Value = 1576.7693
$ws.Range("I85").Value = 1183
$ws.Range("J85").Value = 3742.5
$ws.Range("K85").Value = 1183
$ws.Range("L85").Value = 3742.5
$ws.Range("M85").Value = 65
$ws.Range("N85").Value = -6238.5
$ws.Range("H100").Value = 4980
$ws.Range("J100").Value = 4980
$ws.Range("L100").Value = 4980
$ws.Range("N100").Value = -6062
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1504.6666
$ws.Range("I96").Value = 985
$ws.Range("J96").Value = 1653.1428
$ws.Range("K96").Value = 985
$ws.Range("L96").Value = 1653.1428
$ws.Range("M96").Value = 388
$ws.Range("N96").Value = -4399.1428
